$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (changed) date column for every existing data row
# (rows 2-441) from 2023-09-17 (45186) to 2023-09-19 (45188).
$ws.Range("C2:C441").Value = 45188

# Row 441 picks up an explicit row height now that it is no longer the
# last row in the sheet.
$ws.Rows.Item(441).RowHeight = 15

# Append the new record as row 442.
$ws.Cells.Item(442, 1).Value = "A 43722-2023"

$ws.Cells.Item(442, 2).Value = 45187
$ws.Cells.Item(442, 2).NumberFormat = "YYYY-MM-DD"

$ws.Cells.Item(442, 3).Value = 45188
$ws.Cells.Item(442, 3).NumberFormat = "YYYY-MM-DD"

$ws.Cells.Item(442, 4).Value = "UPPSALA LÄN"
$ws.Cells.Item(442, 5).Value = "TIERP"
$ws.Cells.Item(442, 6).Value = "Bergvik skog öst AB"

$ws.Cells.Item(442, 7).Value = 26.7
$ws.Cells.Item(442, 8).Value = 0
$ws.Cells.Item(442, 9).Value = 0
$ws.Cells.Item(442, 10).Value = 0
$ws.Cells.Item(442, 11).Value = 0
$ws.Cells.Item(442, 12).Value = 0
$ws.Cells.Item(442, 13).Value = 0
$ws.Cells.Item(442, 14).Value = 0
$ws.Cells.Item(442, 15).Value = 0
$ws.Cells.Item(442, 16).Value = 0
$ws.Cells.Item(442, 17).Value = 0

# Column R keeps the wrap-text style applied across the whole table, even
# though there is no "Artnamn" text for this row.
$ws.Cells.Item(442, 18).Value = ""
$ws.Cells.Item(442, 18).WrapText = $true
